# Calibrated bldgs/SoCEUtiNTY and trans/SoCDTtiNTY
# Converts the formula-linked "Share that is New" values on the
# SoCDTtiNTY-psgr / SoCDTtiNTY-frgt sheets into their newly-calibrated,
# hard-coded numbers (breaking the link to 'Calibration Helper'!B70:H75 /
# 'Calibration Helper'!B80:H85).

$wb = $excel.ActiveWorkbook

$psgr = $wb.Worksheets.Item("SoCDTtiNTY-psgr")
$frgt = $wb.Worksheets.Item("SoCDTtiNTY-frgt")
$calib = $wb.Worksheets.Item("Calibration Helper")
$about = $wb.Worksheets.Item("About")

# --- SoCDTtiNTY-psgr (passenger) : rows 2-7, cols B-H ---
$psgr.Range("B2:H2").Value = 0.091875556
$psgr.Range("B3:H3").Value = 0.046942935
$psgr.Range("B4:H4").Value = 0.069375556
$psgr.Range("B5:H5").Value = 0.029411765
$psgr.Range("B6:H6").Value = 0.03030303
$psgr.Range("B7:H7").Value = 0.070625556

# --- SoCDTtiNTY-frgt (freight) : rows 2-7, cols B-H ---
$frgt.Range("B2:H2").Value = 0.085555556
$frgt.Range("B3:H3").Value = 0.035714286
$frgt.Range("B4:H4").Value = 0.065555556
$frgt.Range("B5:H5").Value = 0.029411766
$frgt.Range("B6:H6").Value = 0.03030303
$frgt.Range("B7:H7").Value = 0

# --- Replicate the author's final view / selection state ---
$about.Activate() | Out-Null
$about.Range("J13").Select() | Out-Null

$calib.Activate() | Out-Null
$calib.Range("C23").Select() | Out-Null

$frgt.Activate() | Out-Null
$frgt.Range("O15").Select() | Out-Null

$psgr.Activate() | Out-Null
$psgr.Range("E14").Select() | Out-Null
